$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Ntn4 -> MuSCs -> Dcc)
$ws.Range("G2").Value = 1.211767666666667
$ws.Range("H2").Value = 3.635303
$ws.Range("I2").Value = 0.0191872416143265
$ws.Range("J2").Value = 0.0191872416143265
$ws.Range("Q2").Value = 0.1600385596592222
$ws.Range("R2").Value = 1.440347036933
$ws.Range("S2").Value = 0.0191872416143265
$ws.Range("T2").Value = 0.0191872416143265

# Row 3
$ws.Range("G3").Value = 27.75404733333334
$ws.Range("H3").Value = 83.26214200000001
$ws.Range("I3").Value = 0.4394601594090953
$ws.Range("J3").Value = 0.4394601594090954
$ws.Range("Q3").Value = 3.665486282662444
$ws.Range("S3").Value = 0.4394601594090953
$ws.Range("T3").Value = 0.4394601594090954

# Row 4
$ws.Range("G4").Value = 34.18905066666667
$ws.Range("H4").Value = 102.567152
$ws.Range("I4").Value = 0.5413525989765782
$ws.Range("J4").Value = 0.5413525989765782
$ws.Range("Q4").Value = 4.515359317896888
$ws.Range("R4").Value = 40.638233861072
$ws.Range("S4").Value = 0.5413525989765782
$ws.Range("T4").Value = 0.5413525989765782
